# Apply the "10 nodes 300 sec 90 prob" results into the "10nodes" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10nodes")

# Fill A3:A102 with the simulation results (mostly zeros, with a couple of
# successful deliveries recorded in A44/B44 and A70/B70).
$ws.Range("A3:A102").Value = 0

$ws.Range("A44").Value = 4
$ws.Range("B44").Value = 297

$ws.Range("A70").Value = 2
$ws.Range("B70").Value = 289

# Move the active selection from the data entry range to the summary row.
$ws.Range("A103").Select() | Out-Null
